# "Proceso 1 done" — the A-column running counter (row 2 = base value,
# rows 3:51 = "=previous+1" shared formula) is shifted from starting at 1
# to starting at 151, so the visible sequence becomes 151..200 instead of
# 1..50. Only the seed cell A2 needs to change; the dependent formulas in
# A3:A51 recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A2").Value = 151

# Leave the cursor/selection on A3, matching the saved sheet view.
[void]$ws.Range("A3").Select()
